$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.926.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.628.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.79%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.87'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.68%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.625.16'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.97%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.56%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.21'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.71'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.23%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.70%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.816.21'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.62%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.17%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '357.06'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.88'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +9.68%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.67%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.08%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '548.62'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.35'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.136'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.73%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.50'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.49%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.74%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.32%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.15'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.72%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.37%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.23'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.22%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.42'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0297'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.579'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '151.47'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.77'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.48%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.28%  '
